$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, shifting the "pairwise" comparison blocks down
[void]$ws.Rows.Item(8).Insert()

# Add a new pairwise-id entry in the newly inserted row, reusing the
# existing "aa" / "aaa" pair already used at C2 / D3
$ws.Range("C8").Value = "aa"
$ws.Range("D8").Value = "aaa"

# Update the active selection to match the new layout
[void]$ws.Range("D10").Select()

Write-Output "done"
